$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two MATLAB toolbox line items (rows 2 and 3) are being consolidated into
# a single "Matlab License" line item. Remove row 3 ("Neural Networks" /
# "MATLAB Neural Network Toolbox") entirely, which shifts every following row
# up by one (ADDITIONS header moves from row 4 to row 3, the remaining
# purchase rows move from 5-8 to 4-7, the Total row moves from 9 to 8, etc.).
$ws.Rows("3:3").Delete()

# Update the remaining first line item (row 2) to reflect the new combined
# Matlab license entry and its cost; the Subtotal formula (D2*C2) recalculates
# automatically.
$ws.Range("A2").Value = "Matlab License"
$ws.Range("B2").Value = "University Full Matlab Add-Ons"
$ws.Range("D2").Value = 250

# Match the author's cursor/selection position after the edit.
$ws.Range("A10").Select()
